# Auto-generated Excel COM-interop script
# Applies the numeric cell-value changes described by the source diff
# (scheduled-runner update to the per-sheet profit/price columns H..N).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4679
$ws.Range("I2").Value = 3200.25
$ws.Range("J2").Value = 6650.6665
$ws.Range("K2").Value = 3200.25
$ws.Range("L2").Value = 6650.6665
$ws.Range("M2").Value = -3087.25
$ws.Range("N2").Value = -6876.6665
$ws.Range("H4").Value = 451.13513
$ws.Range("I4").Value = 240.37038
$ws.Range("K4").Value = 240.37038
$ws.Range("M4").Value = -126.37038
$ws.Range("H17").Value = 3041
$ws.Range("I17").Value = 2000
$ws.Range("J17").Value = 3115.3572
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 9346.071599999999
$ws.Range("M17").Value = -5832
$ws.Range("N17").Value = -9682.071599999999
$ws.Range("H19").Value = 7212.8125
$ws.Range("I19").Value = 3872.0952
$ws.Range("J19").Value = 13590.546
$ws.Range("K19").Value = 3872.0952
$ws.Range("L19").Value = 13590.546
$ws.Range("M19").Value = -3697.0952
$ws.Range("N19").Value = -13940.546
$ws.Range("H51").Value = 26742.215
$ws.Range("J51").Value = 31563.182
$ws.Range("L51").Value = 31563.182
$ws.Range("N51").Value = -32531.182
$ws.Range("H53").Value = 408.91666
$ws.Range("I53").Value = 258.64285
$ws.Range("K53").Value = 258.64285
$ws.Range("M53").Value = 378.35715
$ws.Range("H64").Value = 3340.4
$ws.Range("I64").Value = 3400.6667
$ws.Range("J64").Value = 3250
$ws.Range("K64").Value = 3400.6667
$ws.Range("L64").Value = 3250
$ws.Range("M64").Value = -3152.6667
$ws.Range("N64").Value = -3746
$ws.Range("H67").Value = 3340.4
$ws.Range("I67").Value = 3400.6667
$ws.Range("J67").Value = 3250
$ws.Range("K67").Value = 3400.6667
$ws.Range("L67").Value = 3250
$ws.Range("M67").Value = -2542.6667
$ws.Range("N67").Value = -4966
$ws.Range("H69").Value = 20209.637
$ws.Range("H72").Value = 20209.637
$ws.Range("H74").Value = 6933.6
$ws.Range("H77").Value = 6933.6
$ws.Range("H86").Value = 14061.875
$ws.Range("I86").Value = 14856.429
$ws.Range("K86").Value = 14856.429
$ws.Range("M86").Value = -13733.429
$ws.Range("H88").Value = 1255.2
$ws.Range("I88").Value = 650.75
$ws.Range("J88").Value = 1658.1666
$ws.Range("K88").Value = 650.75
$ws.Range("L88").Value = 1658.1666
$ws.Range("M88").Value = -244.75
$ws.Range("N88").Value = -2470.1666
$ws.Range("H89").Value = 14061.875
$ws.Range("I89").Value = 14856.429
$ws.Range("K89").Value = 74282.145
$ws.Range("M89").Value = -68666.145
$ws.Range("H91").Value = 1255.2
$ws.Range("I91").Value = 650.75
$ws.Range("J91").Value = 1658.1666
$ws.Range("K91").Value = 650.75
$ws.Range("L91").Value = 1658.1666
$ws.Range("M91").Value = 753.25
$ws.Range("N91").Value = -4466.1666
$ws.Range("H98").Value = 6869.4443
$ws.Range("I98").Value = 7642.7085
$ws.Range("K98").Value = 7642.7085
$ws.Range("M98").Value = -6144.7085
$ws.Range("H112").Value = 2722.72
$ws.Range("J112").Value = 2722.72
$ws.Range("L112").Value = 8168.16
$ws.Range("N112").Value = -10384.16
$ws.Range("H122").Value = 6869.4443
$ws.Range("I122").Value = 7642.7085
$ws.Range("K122").Value = 22928.1255
$ws.Range("M122").Value = -20478.1255
$ws.Range("H123").Value = 124991.75
$ws.Range("J123").Value = 124991.75
$ws.Range("L123").Value = 124991.75
$ws.Range("N123").Value = -134791.75
$ws.Range("H125").Value = 2755.4285
$ws.Range("I125").Value = 2832.111
$ws.Range("K125").Value = 25488.999
$ws.Range("M125").Value = -23028.999
$ws.Range("H128").Value = 50779.953
$ws.Range("J128").Value = 50779.953
$ws.Range("L128").Value = 50779.953
$ws.Range("N128").Value = -60739.953
$ws.Range("H137").Value = 5561.4614
$ws.Range("I137").Value = 4499.0835
$ws.Range("J137").Value = 6472.0713
$ws.Range("K137").Value = 13497.2505
$ws.Range("L137").Value = 19416.2139
$ws.Range("M137").Value = -10947.2505
$ws.Range("N137").Value = -24516.2139
$ws.Range("H138").Value = 6026.6445
$ws.Range("I138").Value = 4756
$ws.Range("J138").Value = 6544.315
$ws.Range("K138").Value = 14268
$ws.Range("L138").Value = 19632.945
$ws.Range("M138").Value = -9128
$ws.Range("N138").Value = -29912.945
$ws.Range("H141").Value = 1708.6129
$ws.Range("I141").Value = 1570.1786
$ws.Range("J141").Value = 3000.6667
$ws.Range("K141").Value = 4710.5358
$ws.Range("L141").Value = 9002.000100000001
$ws.Range("M141").Value = 469.4642000000003
$ws.Range("N141").Value = -19362.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 46586.285
$ws.Range("I32").Value = 46791.066
$ws.Range("K32").Value = 46791.066
$ws.Range("M32").Value = -46504.066
$ws.Range("H45").Value = 1709.1923
$ws.Range("I45").Value = 1620
$ws.Range("J45").Value = 1909.875
$ws.Range("K45").Value = 1620
$ws.Range("L45").Value = 1909.875
$ws.Range("M45").Value = -1243
$ws.Range("N45").Value = -2663.875
$ws.Range("H61").Value = 9057.454
$ws.Range("I61").Value = 4819.7
$ws.Range("J61").Value = 12588.917
$ws.Range("K61").Value = 4819.7
$ws.Range("L61").Value = 12588.917
$ws.Range("M61").Value = -4607.7
$ws.Range("N61").Value = -13012.917
$ws.Range("H74").Value = 418862.03
$ws.Range("I74").Value = 589764.5
$ws.Range("K74").Value = 589764.5
$ws.Range("M74").Value = -588890.5
$ws.Range("H77").Value = 418862.03
$ws.Range("I77").Value = 589764.5
$ws.Range("K77").Value = 2948822.5
$ws.Range("M77").Value = -2944454.5
$ws.Range("H80").Value = 50000
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996
$ws.Range("H83").Value = 50000
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984
$ws.Range("H97").Value = 1308952.6
$ws.Range("I97").Value = 1852345.8
$ws.Range("K97").Value = 1852345.8
$ws.Range("M97").Value = -1851849.8
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2242.2896
$ws.Range("I122").Value = 1442.9131
$ws.Range("K122").Value = 4328.7393
$ws.Range("M122").Value = -1878.7393
$ws.Range("H132").Value = 9522.558000000001
$ws.Range("I132").Value = 4126.511
$ws.Range("K132").Value = 12379.533
$ws.Range("M132").Value = -9849.533000000001
$ws.Range("H136").Value = 9057.454
$ws.Range("I136").Value = 4819.7
$ws.Range("J136").Value = 12588.917
$ws.Range("K136").Value = 14459.1
$ws.Range("L136").Value = 37766.751
$ws.Range("M136").Value = -11909.1
$ws.Range("N136").Value = -42866.751

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2570.7144
$ws.Range("I22").Value = 2666
$ws.Range("K22").Value = 2666
$ws.Range("M22").Value = -2493
$ws.Range("H24").Value = 1011.25
$ws.Range("I24").Value = 1015.6667
$ws.Range("J24").Value = 998
$ws.Range("K24").Value = 1015.6667
$ws.Range("L24").Value = 998
$ws.Range("M24").Value = -780.6667
$ws.Range("N24").Value = -1468
$ws.Range("H25").Value = 543
$ws.Range("I25").Value = 90
$ws.Range("K25").Value = 90
$ws.Range("M25").Value = 145
$ws.Range("H99").Value = 32273.625
$ws.Range("I99").Value = 32273.625
$ws.Range("K99").Value = 32273.625
$ws.Range("M99").Value = -30775.625
$ws.Range("H105").Value = 29421714
$ws.Range("I105").Value = 41679130
$ws.Range("K105").Value = 41679130
$ws.Range("M105").Value = -41677383
$ws.Range("H107").Value = 1895.6389
$ws.Range("I107").Value = 1431.6666
$ws.Range("K107").Value = 1431.6666
$ws.Range("M107").Value = 488.3334
$ws.Range("H134").Value = 6986.2896
$ws.Range("I134").Value = 6056.269
$ws.Range("K134").Value = 18168.807
$ws.Range("M134").Value = -15633.807

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8498
$ws.Range("J4").Value = 9663.666999999999
$ws.Range("L4").Value = 9663.666999999999
$ws.Range("N4").Value = -9887.666999999999
$ws.Range("H7").Value = 481.81818
$ws.Range("I7").Value = 535.05554
$ws.Range("K7").Value = 535.05554
$ws.Range("M7").Value = -422.05554
$ws.Range("H13").Value = 2864.0256
$ws.Range("I13").Value = 3086.6287
$ws.Range("J13").Value = 916.25
$ws.Range("K13").Value = 3086.6287
$ws.Range("L13").Value = 916.25
$ws.Range("M13").Value = -2947.6287
$ws.Range("N13").Value = -1194.25
$ws.Range("H14").Value = 3375.2727
$ws.Range("J14").Value = 2104.6667
$ws.Range("L14").Value = 2104.6667
$ws.Range("N14").Value = -2444.6667
$ws.Range("H16").Value = 2577.0967
$ws.Range("I16").Value = 2107.8262
$ws.Range("K16").Value = 2107.8262
$ws.Range("M16").Value = -1820.8262
$ws.Range("H28").Value = 50000
$ws.Range("J28").Value = 50000
$ws.Range("L28").Value = 50000
$ws.Range("N28").Value = -50490
$ws.Range("H31").Value = 6884.85
$ws.Range("J31").Value = 7770.4287
$ws.Range("L31").Value = 7770.4287
$ws.Range("N31").Value = -8360.4287
$ws.Range("H34").Value = 6884.85
$ws.Range("J34").Value = 7770.4287
$ws.Range("L34").Value = 7770.4287
$ws.Range("N34").Value = -8174.4287
$ws.Range("H58").Value = 5717.6113
$ws.Range("I58").Value = 4280.364
$ws.Range("J58").Value = 7976.143
$ws.Range("K58").Value = 4280.364
$ws.Range("L58").Value = 7976.143
$ws.Range("M58").Value = -4077.364
$ws.Range("N58").Value = -8382.143
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H99").Value = 5199.077
$ws.Range("I99").Value = 6696
$ws.Range("K99").Value = 6696
$ws.Range("M99").Value = -5198
$ws.Range("H105").Value = 736.9474
$ws.Range("J105").Value = 907
$ws.Range("L105").Value = 907
$ws.Range("N105").Value = -4401
$ws.Range("H107").Value = 498
$ws.Range("I107").Value = 490
$ws.Range("J107").Value = 514
$ws.Range("K107").Value = 490
$ws.Range("L107").Value = 514
$ws.Range("M107").Value = 1430
$ws.Range("N107").Value = -4354
$ws.Range("H113").Value = 2577.0967
$ws.Range("I113").Value = 2107.8262
$ws.Range("K113").Value = 2107.8262
$ws.Range("M113").Value = 62.17380000000003
$ws.Range("H122").Value = 1096.7693
$ws.Range("I122").Value = 945.7778
$ws.Range("K122").Value = 2837.3334
$ws.Range("M122").Value = -387.3334
$ws.Range("H126").Value = 5199.077
$ws.Range("I126").Value = 6696
$ws.Range("K126").Value = 20088
$ws.Range("M126").Value = -17618
$ws.Range("H132").Value = 20225.824
$ws.Range("I132").Value = 3499.3928
$ws.Range("J132").Value = 98282.5
$ws.Range("K132").Value = 10498.1784
$ws.Range("L132").Value = 294847.5
$ws.Range("M132").Value = -7968.178400000001
$ws.Range("N132").Value = -299907.5
$ws.Range("H136").Value = 5717.6113
$ws.Range("I136").Value = 4280.364
$ws.Range("J136").Value = 7976.143
$ws.Range("K136").Value = 12841.092
$ws.Range("L136").Value = 23928.429
$ws.Range("M136").Value = -10291.092
$ws.Range("N136").Value = -29028.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4704079
$ws.Range("I4").Value = 3335641.2
$ws.Range("J4").Value = 7745051.5
$ws.Range("K4").Value = 10006923.6
$ws.Range("L4").Value = 23235154.5
$ws.Range("M4").Value = -10006811.6
$ws.Range("N4").Value = -23235378.5
$ws.Range("H5").Value = 464
$ws.Range("I5").Value = 299.53333
$ws.Range("K5").Value = 898.5999899999999
$ws.Range("M5").Value = -786.5999899999999
$ws.Range("H11").Value = 163416.23
$ws.Range("I11").Value = 12434.777
$ws.Range("K11").Value = 37304.331
$ws.Range("M11").Value = -37164.331
$ws.Range("H112").Value = 6673341.5
$ws.Range("I112").Value = 6673341.5
$ws.Range("K112").Value = 20020024.5
$ws.Range("M112").Value = -20018916.5
$ws.Range("H129").Value = 62501484
$ws.Range("I129").Value = 841.75
$ws.Range("J129").Value = 125002130
$ws.Range("K129").Value = 2525.25
$ws.Range("L129").Value = 375006390
$ws.Range("M129").Value = 2474.75
$ws.Range("N129").Value = -375016390
$ws.Range("H132").Value = 35309.875
$ws.Range("I132").Value = 78139.78999999999
$ws.Range("J132").Value = 1997.7222
$ws.Range("K132").Value = 703258.11
$ws.Range("L132").Value = 17979.4998
$ws.Range("M132").Value = -700728.11
$ws.Range("N132").Value = -23039.4998
$ws.Range("H135").Value = 464
$ws.Range("I135").Value = 299.53333
$ws.Range("K135").Value = 2695.79997
$ws.Range("M135").Value = -160.79997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 144230.28
$ws.Range("I2").Value = 1602.1666
$ws.Range("J2").Value = 999999
$ws.Range("K2").Value = 1602.1666
$ws.Range("L2").Value = 999999
$ws.Range("M2").Value = -1489.1666
$ws.Range("N2").Value = -1000225
$ws.Range("H4").Value = 15
$ws.Range("I4").Value = 15
$ws.Range("K4").Value = 15
$ws.Range("M4").Value = 97
$ws.Range("H21").Value = 7569
$ws.Range("J21").Value = 7569
$ws.Range("L21").Value = 7569
$ws.Range("N21").Value = -7915
$ws.Range("H30").Value = 7569
$ws.Range("J30").Value = 7569
$ws.Range("L30").Value = 7569
$ws.Range("N30").Value = -7779
$ws.Range("H80").Value = 4002.5
$ws.Range("I80").Value = 6005
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 6005
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -5007
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 4002.5
$ws.Range("I83").Value = 6005
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 30025
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -25033
$ws.Range("N83").Value = -19984
$ws.Range("H97").Value = 4059.6667
$ws.Range("I97").Value = 2181.8
$ws.Range("K97").Value = 2181.8
$ws.Range("M97").Value = -1685.8
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 3722.25
$ws.Range("I113").Value = 2449.5
$ws.Range("J113").Value = 4995
$ws.Range("K113").Value = 2449.5
$ws.Range("L113").Value = 4995
$ws.Range("M113").Value = -279.5
$ws.Range("N113").Value = -9335
$ws.Range("H116").Value = 72000
$ws.Range("J116").Value = 72000
$ws.Range("L116").Value = 72000
$ws.Range("N116").Value = -81178
$ws.Range("H122").Value = 2397.1
$ws.Range("I122").Value = 2202.6428
$ws.Range("J122").Value = 2850.8333
$ws.Range("K122").Value = 6607.928400000001
$ws.Range("L122").Value = 8552.499899999999
$ws.Range("M122").Value = -4157.928400000001
$ws.Range("N122").Value = -13452.4999
$ws.Range("H132").Value = 4985.3184
$ws.Range("I132").Value = 3246.9312
$ws.Range("J132").Value = 8346.200000000001
$ws.Range("K132").Value = 9740.793600000001
$ws.Range("L132").Value = 25038.6
$ws.Range("M132").Value = -7210.793600000001
$ws.Range("N132").Value = -30098.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4279.1113
$ws.Range("I7").Value = 4237.4287
$ws.Range("K7").Value = 4237.4287
$ws.Range("M7").Value = -4125.4287
$ws.Range("H16").Value = 987.5278
$ws.Range("I16").Value = 948.55884
$ws.Range("K16").Value = 948.55884
$ws.Range("M16").Value = -778.55884
$ws.Range("H22").Value = 3506.9395
$ws.Range("I22").Value = 2184.25
$ws.Range("J22").Value = 4262.7617
$ws.Range("K22").Value = 2184.25
$ws.Range("L22").Value = 4262.7617
$ws.Range("M22").Value = -1889.25
$ws.Range("N22").Value = -4852.7617
$ws.Range("H27").Value = 3506.9395
$ws.Range("I27").Value = 2184.25
$ws.Range("J27").Value = 4262.7617
$ws.Range("K27").Value = 2184.25
$ws.Range("L27").Value = 4262.7617
$ws.Range("M27").Value = -2077.25
$ws.Range("N27").Value = -4476.7617
$ws.Range("H40").Value = 33337846
$ws.Range("I40").Value = 35718828
$ws.Range("K40").Value = 35718828
$ws.Range("M40").Value = -35718692
$ws.Range("H46").Value = 4466.4287
$ws.Range("J46").Value = 7700.1
$ws.Range("L46").Value = 7700.1
$ws.Range("N46").Value = -8076.1
$ws.Range("H55").Value = 712.1875
$ws.Range("I55").Value = 709.4
$ws.Range("J55").Value = 716.8333
$ws.Range("K55").Value = 709.4
$ws.Range("L55").Value = 716.8333
$ws.Range("M55").Value = -536.4
$ws.Range("N55").Value = -1062.8333
$ws.Range("H68").Value = 2997.1428
$ws.Range("J68").Value = 2990
$ws.Range("L68").Value = 2990
$ws.Range("N68").Value = -4488
$ws.Range("H71").Value = 2997.1428
$ws.Range("J71").Value = 2990
$ws.Range("L71").Value = 14950
$ws.Range("N71").Value = -22438
$ws.Range("H82").Value = 1642.1428
$ws.Range("I82").Value = 1665.8334
$ws.Range("K82").Value = 1665.8334
$ws.Range("M82").Value = -1304.8334
$ws.Range("H85").Value = 1642.1428
$ws.Range("I85").Value = 1665.8334
$ws.Range("K85").Value = 1665.8334
$ws.Range("M85").Value = -417.8334
$ws.Range("H122").Value = 45458944
$ws.Range("I122").Value = 62504188
$ws.Range("J122").Value = 4966.6665
$ws.Range("K122").Value = 187512564
$ws.Range("L122").Value = 14899.9995
$ws.Range("M122").Value = -187510114
$ws.Range("N122").Value = -19799.9995
$ws.Range("H126").Value = 4279.1113
$ws.Range("I126").Value = 4237.4287
$ws.Range("K126").Value = 12712.2861
$ws.Range("M126").Value = -10242.2861
$ws.Range("H132").Value = 7363.6343
$ws.Range("I132").Value = 6590.2964
$ws.Range("K132").Value = 19770.8892
$ws.Range("M132").Value = -17240.8892
$ws.Range("H136").Value = 9214.625
$ws.Range("I136").Value = 8267.362999999999
$ws.Range("K136").Value = 24802.089
$ws.Range("L136").Value = 33895.8
$ws.Range("M136").Value = -22252.089
$ws.Range("N136").Value = -38995.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 4994
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H54").Value = 94443.89
$ws.Range("I54").Value = 54666.668
$ws.Range("J54").Value = 114332.5
$ws.Range("K54").Value = 54666.668
$ws.Range("L54").Value = 114332.5
$ws.Range("M54").Value = -54146.668
$ws.Range("N54").Value = -115372.5
$ws.Range("H81").Value = 9744.032999999999
$ws.Range("I81").Value = 6954.8
$ws.Range("K81").Value = 13909.6
$ws.Range("M81").Value = -12848.6
$ws.Range("H84").Value = 9744.032999999999
$ws.Range("I84").Value = 6954.8
$ws.Range("K84").Value = 69548
$ws.Range("M84").Value = -64244
$ws.Range("H107").Value = 1891.6154
$ws.Range("I107").Value = 1838.5714
$ws.Range("K107").Value = 5515.7142
$ws.Range("M107").Value = -3595.7142
$ws.Range("H132").Value = 4469.0894
$ws.Range("I132").Value = 4048.5
$ws.Range("K132").Value = 12145.5
$ws.Range("M132").Value = -9615.5
$ws.Range("H136").Value = 3558.4666
$ws.Range("I136").Value = 2774.1904
$ws.Range("J136").Value = 5388.4443
$ws.Range("K136").Value = 8322.5712
$ws.Range("L136").Value = 16165.3329
$ws.Range("M136").Value = -5772.5712
$ws.Range("N136").Value = -21265.3329

Write-Output "Applied 500 value updates and 4 clears"
